$p = $ppt.ActivePresentation

# Slide 9 ("Página adicionar aula" / "Página alterar aula" mockups, currently
# mislabeled as the Professor's home page) -> becomes the admin page for
# adding/editing classes.
$s9 = $p.Slides.Item(9)
$sh9 = $s9.Shapes.Item(8)
$sh9.Width = 434.6508661417323
$sh9.TextFrame.TextRange.Text = "Página adicionar e editar aulas – Visão Administrador"

# Slide 10 ("controle de aulas / presenças" mockup) -> becomes the page to
# view and manually edit attendance.
$s10 = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(9)
$sh10.Width = 595.4782677165355
$sh10.TextFrame.TextRange.Text = "Página para visualizar e alterar presenças manualmente – Visão Professor"
